# Framework changes made to search module
#
# "B Suite" (row 3) is renamed to "Search" and "E Suite" (row 6) is renamed
# to "Watchlist" on the "Test Suite" sheet, to align the TSID column with
# the module names already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update A6 first, then A3, so new shared-string entries are appended to
# the workbook's string table in the same order produced by the original
# edit (Watchlist, then Search).
$ws.Range("A6").Value = "Watchlist"
$ws.Range("A3").Value = "Search"

# Reflect the saved window size from the edited workbook (best effort —
# harmless if the host doesn't expose window geometry).
try {
    $excel.ActiveWindow.Width = 10305
    $excel.ActiveWindow.Height = 3150
} catch {}
